$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Igfbp4"
$ws.Range("C2").Value = "Fzd8"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 9.568356333333332
$ws.Range("H2").Value = 28.705069
$ws.Range("I2").Value = 0.09973288675158326
$ws.Range("J2").Value = 0.09973288675158326
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 2.558821666666667
$ws.Range("N2").Value = 7.676465
$ws.Range("O2").Value = 0.2156728774407755
$ws.Range("P2").Value = 0.2156728774407755
$ws.Range("Q2").Value = 24.48371750012056
$ws.Range("R2").Value = 220.353457501085
$ws.Range("S2").Value = 0.02150967866118896
$ws.Range("T2").Value = 0.02150967866118896

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Igfbp4"
$ws.Range("C3").Value = "Fzd8"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 9.568356333333332
$ws.Range("H3").Value = 28.705069
$ws.Range("I3").Value = 0.09973288675158326
$ws.Range("J3").Value = 0.09973288675158326
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 6.453984666666667
$ws.Range("N3").Value = 19.361954
$ws.Range("O3").Value = 0.5439806384912759
$ws.Range("P3").Value = 0.5439806384912759
$ws.Range("Q3").Value = 61.75402506053621
$ws.Range("R3").Value = 555.786225544826
$ws.Range("S3").Value = 0.05425275941370437
$ws.Range("T3").Value = 0.05425275941370437

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Igfbp4"
$ws.Range("C4").Value = "Fzd8"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 9.568356333333332
$ws.Range("H4").Value = 28.705069
$ws.Range("I4").Value = 0.09973288675158326
$ws.Range("J4").Value = 0.09973288675158326
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 2.851558333333334
$ws.Range("N4").Value = 8.554675000000001
$ws.Range("O4").Value = 0.2403464840679487
$ws.Range("P4").Value = 0.2403464840679487
$ws.Range("Q4").Value = 27.28472623861945
$ws.Range("R4").Value = 245.562536147575
$ws.Range("S4").Value = 0.02397044867668994
$ws.Range("T4").Value = 0.02397044867668994

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Igfbp4"
$ws.Range("C5").Value = "Fzd8"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 75.94550566666668
$ws.Range("H5").Value = 227.836517
$ws.Range("I5").Value = 0.7915951551217724
$ws.Range("J5").Value = 0.7915951551217723
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 2.558821666666667
$ws.Range("N5").Value = 7.676465
$ws.Range("O5").Value = 0.2156728774407755
$ws.Range("P5").Value = 0.2156728774407755
$ws.Range("Q5").Value = 194.3310053858228
$ws.Range("R5").Value = 1748.979048472405
$ws.Range("S5").Value = 0.1707256048732896
$ws.Range("T5").Value = 0.1707256048732896

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Igfbp4"
$ws.Range("C6").Value = "Fzd8"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 75.94550566666668
$ws.Range("H6").Value = 227.836517
$ws.Range("I6").Value = 0.7915951551217724
$ws.Range("J6").Value = 0.7915951551217723
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 6.453984666666667
$ws.Range("N6").Value = 19.361954
$ws.Range("O6").Value = 0.5439806384912759
$ws.Range("P6").Value = 0.5439806384912759
$ws.Range("Q6").Value = 490.1511290749132
$ws.Range("R6").Value = 4411.360161674218
$ws.Range("S6").Value = 0.4306124379097423
$ws.Range("T6").Value = 0.4306124379097422

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Igfbp4"
$ws.Range("C7").Value = "Fzd8"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 75.94550566666668
$ws.Range("H7").Value = 227.836517
$ws.Range("I7").Value = 0.7915951551217724
$ws.Range("J7").Value = 0.7915951551217723
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 2.851558333333334
$ws.Range("N7").Value = 8.554675000000001
$ws.Range("O7").Value = 0.2403464840679487
$ws.Range("P7").Value = 0.2403464840679487
$ws.Range("Q7").Value = 216.5630395629973
$ws.Range("R7").Value = 1949.067356066975
$ws.Range("S7").Value = 0.1902571123387404
$ws.Range("T7").Value = 0.1902571123387404

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Igfbp4"
$ws.Range("C8").Value = "Fzd8"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 10.42596933333333
$ws.Range("H8").Value = 31.277908
$ws.Range("I8").Value = 0.1086719581266445
$ws.Range("J8").Value = 0.1086719581266445
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 2.558821666666667
$ws.Range("N8").Value = 7.676465
$ws.Range("O8").Value = 0.2156728774407755
$ws.Range("P8").Value = 0.2156728774407755
$ws.Range("Q8").Value = 26.67819622613555
$ws.Range("R8").Value = 240.10376603522
$ws.Range("S8").Value = 0.02343759390629687
$ws.Range("T8").Value = 0.02343759390629687

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Igfbp4"
$ws.Range("C9").Value = "Fzd8"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 10.42596933333333
$ws.Range("H9").Value = 31.277908
$ws.Range("I9").Value = 0.1086719581266445
$ws.Range("J9").Value = 0.1086719581266445
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 6.453984666666667
$ws.Range("N9").Value = 19.361954
$ws.Range("O9").Value = 0.5439806384912759
$ws.Range("P9").Value = 0.5439806384912759
$ws.Range("Q9").Value = 67.28904621247023
$ws.Range("R9").Value = 605.601415912232
$ws.Range("S9").Value = 0.05911544116782926
$ws.Range("T9").Value = 0.05911544116782926

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Igfbp4"
$ws.Range("C10").Value = "Fzd8"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 10.42596933333333
$ws.Range("H10").Value = 31.277908
$ws.Range("I10").Value = 0.1086719581266445
$ws.Range("J10").Value = 0.1086719581266445
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 2.851558333333334
$ws.Range("N10").Value = 8.554675000000001
$ws.Range("O10").Value = 0.2403464840679487
$ws.Range("P10").Value = 0.2403464840679487
$ws.Range("Q10").Value = 29.73025973554445
$ws.Range("R10").Value = 267.5723376199001
$ws.Range("S10").Value = 0.02611892305251834
$ws.Range("T10").Value = 0.02611892305251835

